$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Insert a new ListParagraph after "After the initial set up ..." with the
#    "Next step will be creating ..." text (runs all tagged lang=en-AU).
# ---------------------------------------------------------------------------
$afterInitial = $d.Paragraphs(3)
$afterInitial.Range.InsertParagraphAfter()
$nextStepPara = $d.Paragraphs(4)
$nextStepRange = $d.Range($nextStepPara.Range.Start, $nextStepPara.Range.End)

$nextStepXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:ind w:firstLineChars="0"/></w:pPr><w:r><w:rPr><w:lang w:val="en-AU"/></w:rPr><w:t xml:space="preserve">Next step will be </w:t></w:r><w:r><w:rPr><w:lang w:val="en-AU"/></w:rPr><w:t>creating</w:t></w:r><w:r><w:rPr><w:lang w:val="en-AU"/></w:rPr><w:t xml:space="preserve"> a basic controller with jumping, and a script that makes the background scroll and repeat. After adjusting the scrolling </w:t></w:r><w:r><w:rPr><w:lang w:val="en-AU"/></w:rPr><w:t>speed,</w:t></w:r><w:r><w:rPr><w:lang w:val="en-AU"/></w:rPr><w:t xml:space="preserve"> I manage to create a parallax effect.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$nextStepRange.InsertXML($nextStepXml)

# ---------------------------------------------------------------------------
# 2) Insert a new ListParagraph after that one describing the jumping bug
#    (includes proofErr grammar/spelling markers from the original capture).
# ---------------------------------------------------------------------------
$nextStepPara2 = $d.Paragraphs(4)
$nextStepPara2.Range.InsertParagraphAfter()
$bugPara = $d.Paragraphs(5)
$bugRange = $d.Range($bugPara.Range.Start, $bugPara.Range.End)

$bugXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:ind w:firstLineChars="0"/></w:pPr><w:r><w:t xml:space="preserve">A problem </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>occur</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve">: the jumping feels wired. </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>So</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> I google and </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>youtube</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> some stuff and by adjusting the gravity in the script I manage to fix it.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$bugRange.InsertXML($bugXml)

# ---------------------------------------------------------------------------
# 3) Add the collectable / scene-management sentence to the final (bookmark)
#    paragraph, inserted before the _GoBack bookmark so the bookmark survives.
# ---------------------------------------------------------------------------
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$insertStart = $lastPara.Range.Start
$insertRange = $d.Range($insertStart, $insertStart)

$addXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>A</w:t></w:r><w:r><w:t>dd prefab and script for obstacle and reward after that and made game manager script and the scene handler script, to spawn reward</w:t></w:r><w:r><w:t>,</w:t></w:r><w:r><w:t xml:space="preserve"> obstacle, </w:t></w:r><w:r><w:t>and scene management.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insertRange.InsertXML($addXml)
